$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YouTube")

$values = @(
    "osCkeNUZgt8",
    "oMXGGmBfkf8",
    "UZw8oRnk0Og",
    "FKmCsysiYRk",
    "Mt239UeAnJg",
    "osCkeNUZgt8",
    "50_jRz8LloI",
    "ztSWLFZ8lFQ",
    "Wb9p6uM57I8",
    "1Uci71lS2LA",
    "-h2Zb5rnHdA",
    "JpOmFQsFNL4",
    "xkTTdGqAHVM",
    "W3et1dgZu6Y",
    "RodEgLRYmhs",
    "jsxc-CqIVqg",
    "064Y5rXGm_s",
    "VZVWeZmlRpM"
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

$templateCell = $ws.Cells.Item($lastRow, 1)

for ($i = 0; $i -lt $values.Count; $i++) {
    $cell = $ws.Cells.Item($startRow + $i, 1)
    $cell.Value = $values[$i]
    $templateCell.Copy()
    $cell.PasteSpecial(-4122)
}
